$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4283-MS-EI-DB-SAR-REC-CTRFD-RNI-FEE+INTEREST-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-1st"
$newShortName = "428r"

# Update product name on both sheets
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update short name (now text instead of number)
$wsInput.Range("B2").Value = $newShortName

# Change selection on input sheet and move active tab to output sheet
$wsInput.Range("B3").Select()
$wsOutput.Activate()
